$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E9: convert inline string "+919510038048" to a numeric value
$ws.Range("E9").Value = 919510038048

# New rows of lead data
$ws.Range("A10").Value = "Yes."
$ws.Range("B10").Value = "Interested"
$ws.Range("C10").Value = "neutral"
$ws.Range("D10").Value = "2025-11-04 20:02:48"
$ws.Range("E10").Value = 919510038048

$ws.Range("A11").Value = "Motor."
$ws.Range("B11").Value = "Interested"
$ws.Range("C11").Value = "neutral"
$ws.Range("D11").Value = "2025-11-05 17:11:40"
$ws.Range("E11").Value = 919510038048

$ws.Range("A12").Value = "3G Nola."
$ws.Range("B12").Value = "Interested"
$ws.Range("C12").Value = "neutral"
$ws.Range("D12").Value = "2025-11-05 17:18:54"
$ws.Range("E12").Value = 919106284482

$ws.Range("A13").Value = "3G Nola."
$ws.Range("B13").Value = "Interested"
$ws.Range("C13").Value = "neutral"
$ws.Range("D13").Value = "2025-11-05 17:57:44"
$ws.Range("E13").Value = 919106284482

$ws.Range("A14").Value = "As."
$ws.Range("B14").Value = "Interested"
$ws.Range("C14").Value = "neutral"
$ws.Range("D14").Value = "2025-11-05 17:57:54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "+919328027733"
